# This script reproduces the "update scripts with new tpm" data refresh for the
# Rtn4 -> Tnfrsf19 ligand-receptor pair sheet: recomputed TPM-based values and an
# added "Target cluster" dimension (ECs and MuSCs) for every sending cluster, which
# grows the data block from 6 rows (r2:r7) to 12 rows (r2:r13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($RowIndex, $Values) {
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $ws.Cells.Item($RowIndex, $i + 1).Value = $Values[$i]
    }
}

Set-Row 2 @("ECs", "Rtn4", "Tnfrsf19", "ECs", 2, 1, 78.25435300000001, 156.508706, 0.1800096451968904, 0.1347002616492669, 1, 0.5, 0.054073, 0.108146, 0.04339331850589854, 0.04339331850589854, 4.231447629769001, 16.925790519076, 0.007811215868162452, 0.005845091356574509)
Set-Row 3 @("ECs", "Rtn4", "Tnfrsf19", "MuSCs", 2, 1, 78.25435300000001, 156.508706, 0.1800096451968904, 0.1347002616492669, 2, 1, 1.1920405, 2.384081, 0.9566066814941014, 0.9566066814941014, 93.28235807729652, 373.1294323091861, 0.1721984293287279, 0.1288551702926924)
Set-Row 4 @("FAPs", "Rtn4", "Tnfrsf19", "ECs", 3, 1, 80.25908033333333, 240.777241, 0.1846211491216604, 0.2072265383236163, 1, 0.5, 0.054073, 0.108146, 0.04339331850589854, 0.04339331850589854, 4.339849250864334, 26.039095505186, 0.0080113243267612, 0.00899224718035147)
Set-Row 5 @("FAPs", "Rtn4", "Tnfrsf19", "MuSCs", 3, 1, 80.25908033333333, 240.777241, 0.1846211491216604, 0.2072265383236163, 2, 1, 1.1920405, 2.384081, 0.9566066814941014, 0.9566066814941014, 95.67207425008684, 574.032445500521, 0.1766098247948992, 0.1982342911432648)
Set-Row 6 @("Inflammatory-Mac", "Rtn4", "Tnfrsf19", "ECs", 3, 1, 73.92583733333333, 221.777512, 0.1700526966117318, 0.1908742948416136, 1, 0.5, 0.054073, 0.108146, 0.04339331850589854, 0.04339331850589854, 3.997391802125334, 23.984350812752, 0.007379150826859812, 0.008282669070650927)
Set-Row 7 @("Inflammatory-Mac", "Rtn4", "Tnfrsf19", "MuSCs", 3, 1, 73.92583733333333, 221.777512, 0.1700526966117318, 0.1908742948416136, 2, 1, 1.1920405, 2.384081, 0.9566066814941014, 0.9566066814941014, 88.12259209774534, 528.7355525864721, 0.162673545784872, 0.1825916257709627)
Set-Row 8 @("MuSCs", "Rtn4", "Tnfrsf19", "ECs", 2, 1, 64.0114765, 128.022953, 0.1472465458029419, 0.1101838083448968, 1, 0.5, 0.054073, 0.108146, 0.04339331850589854, 0.04339331850589854, 3.4612925687845, 13.845170275138, 0.006389516260920436, 0.004781241089702988)
Set-Row 9 @("MuSCs", "Rtn4", "Tnfrsf19", "MuSCs", 2, 1, 64.0114765, 128.022953, 0.1472465458029419, 0.1101838083448968, 2, 1, 1.1920405, 2.384081, 0.9566066814941014, 0.9566066814941014, 76.30427245279826, 305.217089811193, 0.1408570295420215, 0.1054025672551938)
Set-Row 10 @("Neutrophils", "Rtn4", "Tnfrsf19", "ECs", 3, 1, 56.48610166666666, 169.458305, 0.1299358148111217, 0.1458454203955994, 1, 0.5, 0.054073, 0.108146, 0.04339331850589854, 0.04339331850589854, 3.054372975421666, 18.32623785253, 0.005638346197422451, 0.006328716779852915)
Set-Row 11 @("Neutrophils", "Rtn4", "Tnfrsf19", "MuSCs", 3, 1, 56.48610166666666, 169.458305, 0.1299358148111217, 0.1458454203955994, 2, 1, 1.1920405, 2.384081, 0.9566066814941014, 0.9566066814941014, 67.33372087378416, 404.002325242705, 0.1242974686136992, 0.1395167036157465)
Set-Row 12 @("Resolving-Mac", "Rtn4", "Tnfrsf19", "ECs", 3, 1, 81.78626233333334, 245.358787, 0.1881341484556537, 0.2111696764450071, 1, 0.5, 0.054073, 0.108146, 0.04339331850589854, 0.04339331850589854, 4.422428563150334, 26.534571378902, 0.008163765025772183, 0.009163353028765734)
Set-Row 13 @("Resolving-Mac", "Rtn4", "Tnfrsf19", "MuSCs", 3, 1, 81.78626233333334, 245.358787, 0.1881341484556537, 0.2111696764450071, 2, 1, 1.1920405, 2.384081, 0.9566066814941014, 0.9566066814941014, 97.49253704495784, 584.9552222697471, 0.1799703834298816, 0.2020063234162414)
